$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9268046617507935
$ws.Range("B1").Value = 1.748464226722717
$ws.Range("C1").Value = 4.163479804992676
$ws.Range("D1").Value = 3.359662055969238
$ws.Range("E1").Value = 0.3774891793727875
